$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: "*" / "Leverancier" / "Leverancier"
$ws.Range("A5").Value = "*"
$ws.Range("B5").Value = "Leverancier"
$ws.Range("C5").Value = "Leverancier"

# Match the formatting used by the row above it (A4:B4 style) for the new A5:B5 cells
$ws.Range("A4:B4").Copy()
$ws.Range("A5:B5").PasteSpecial(-4122)

# PasteSpecial(Formats) can bring along the old text, so re-assert the values
$ws.Range("A5").Value = "*"
$ws.Range("B5").Value = "Leverancier"

# Column A widened (auto-fit) to accommodate its contents
$ws.Columns.Item(1).EntireColumn.AutoFit()

# Leave active cell on A2, as in the saved file
$ws.Range("A2").Select()
